$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Preguntas:" paragraph - extend the dotted underline so it covers
#    "Id-Test, Id-Categoría-p" (previously only "Id-Categoría-p" carried
#    the dotted underline). The text itself does not change, only the
#    character formatting (and, incidentally, the run boundaries).
# ---------------------------------------------------------------------
$anchor1 = $d.Content
$null = $anchor1.Find.Execute("Id-Test, Id-Categoría-p", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base1 = $anchor1.Start

$seg1a = $d.Range($base1 + 0, $base1 + 5)      # "Id-Te"
$seg1a.Font.Underline = 4                       # wdUnderlineDotted

$seg1b = $d.Range($base1 + 5, $base1 + 6)      # "s"
$seg1b.Font.Underline = 4

$seg1c = $d.Range($base1 + 6, $base1 + 21)     # "t, Id-Categoría"
$seg1c.Font.Underline = 4

$seg1d = $d.Range($base1 + 21, $base1 + 23)    # "-p"
$seg1d.Font.Underline = 4

# ---------------------------------------------------------------------
# 2) "Categorías:" paragraph - add a single underline across
#    " Id-categoría," and " Id-Categoría-Text, Id-Categoría-pregunta".
# ---------------------------------------------------------------------
$anchor2 = $d.Content
$null = $anchor2.Find.Execute("Id-categoría, Id-Categoría-Text, Id-Categoría-pregunta", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base2 = $anchor2.Start - 1                     # back up over the leading space

$seg2a = $d.Range($base2 + 0, $base2 + 14)     # " Id-categoría,"
$seg2a.Font.Underline = 1                       # wdUnderlineSingle

$seg2b = $d.Range($base2 + 14, $base2 + 55)    # " Id-Categoría-Text, Id-Categoría-pregunta"
$seg2b.Font.Underline = 1
